# Update the speaker-credits cell on the sponsor workshop row (I4).
# The text is simplified: the trailing "; " and line-break before
# "Livia Puljak, Damir Važanić" are removed so the names read as one
# continuous, comma-separated list.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I4").Value = "Adriano Friganović, Kata Ivanišević, Livia Puljak, Damir Važanić"

# Reflect the author's final view/selection state on the sheet:
# the window had scrolled so column D is the left-most visible column
# (the sheet view's topLeftCell moved from D1 to D4), and the active
# cell/selection moved from K6 to I4 (the cell that was just edited).
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 4
$ws.Range("I4").Select() | Out-Null
